$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 279, shifting existing row 279 (and everything
# below it) down by one. Excel's row-insert carries formatting of the row
# above down into the new row automatically (matches the s="2" date-style
# cell seen in the target diff).
$ws.Rows.Item(279).EntireRow.Insert()

# Populate the newly inserted row 279 with the new weekly price entry.
$ws.Cells.Item(279, 1).Value = 11
$ws.Cells.Item(279, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(279, 3).Value = "Bíobío"
$ws.Cells.Item(279, 4).Value = 44875
$ws.Cells.Item(279, 5).Value = 8
$ws.Cells.Item(279, 6).Value = 100114014
$ws.Cells.Item(279, 7).Value = "Betarraga"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 1100
$ws.Cells.Item(279, 11).Value = 650
$ws.Cells.Item(279, 12).Value = 700
$ws.Cells.Item(279, 13).Value = 677
$ws.Cells.Item(279, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(279, 15).Value = "Región Metropolitana"
$ws.Cells.Item(279, 16).Value = 135
$ws.Cells.Item(279, 17).Value = 5
$ws.Cells.Item(279, 18).Value = "Hortaliza"
